$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5, 4, 3 (bottom-up) so only the header and the first data row remain
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Update row 2 with the new values
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Cd84"
$ws.Range("C2").Value = "Cd84"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.01702333333333333
$ws.Range("H2").Value = 0.05107
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.01702333333333333
$ws.Range("N2").Value = 0.05107
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0002897938777777778
$ws.Range("R2").Value = 0.0026081449
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
